# "Generate Report for Archive"
#
# 1) Status text "Ready for handoff" -> "In Translation" everywhere it
#    appears (shared string is reused across the Overview sheet's zh-cn /
#    de-de status columns and each language sheet's own Status column).
# 2) Narrow the "Status" columns (Overview!E:F, zh-cn!C, de-de!C) from their
#    old autofit width down to the new narrower autofit width.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Update the status text wherever it appears ---
if ($overview.Range("E2").Value() -eq $oldStatus) { $overview.Range("E2").Value = $newStatus }
if ($overview.Range("F2").Value() -eq $oldStatus) { $overview.Range("F2").Value = $newStatus }
if ($zhcn.Range("C2").Value() -eq $oldStatus)      { $zhcn.Range("C2").Value = $newStatus }
if ($dede.Range("C2").Value() -eq $oldStatus)      { $dede.Range("C2").Value = $newStatus }

# --- Narrow the Status columns to match the new (shorter) text ---
$overview.Columns.Item(5).ColumnWidth = 12.5   # Overview!E (zh-cn status)
$overview.Columns.Item(6).ColumnWidth = 12.5   # Overview!F (de-de status)
$zhcn.Columns.Item(3).ColumnWidth = 12.5        # zh-cn!C (Status)
$dede.Columns.Item(3).ColumnWidth = 12.5        # de-de!C (Status)
